$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "segnen",
    "betteln",
    "fesseln",
    "spielen",
    "liegen",
    "stehlen",
    "weichen",
    "kennen",
    "ehren",
    "wenden",
    "bergen",
    "filmen",
    "wohnen",
    "bitten",
    "töten",
    "suchen",
    "drehen",
    "stimmen",
    "jagen",
    "kümmern",
    "ändern",
    "wirken",
    "klettern",
    "arten",
    "parken",
    "schließen",
    "heilen",
    "heißen",
    "ächzen",
    "quälen",
    "bluten",
    "nennen",
    "sprechen",
    "strahlen",
    "bauen",
    "spüren",
    "starten",
    "seufzen",
    "pfeifen",
    "liefern",
    "formen",
    "wachsen",
    "stören",
    "wehtun",
    "machen",
    "reizen",
    "gründen",
    "enden",
    "zielen",
    "gleichen",
    "klingen",
    "landen",
    "schwingen",
    "flüchten",
    "passen",
    "lesen",
    "melden",
    "streichen",
    "platzen",
    "nähen",
    "freuen",
    "decken",
    "zahlen",
    "hupen",
    "weigern",
    "drücken",
    "loben",
    "trauen",
    "hauen",
    "pflanzen",
    "hören",
    "schneiden",
    "achten",
    "locken",
    "folgen",
    "schlucken",
    "sorgen",
    "tropfen",
    "schätzen",
    "heulen",
    "hassen",
    "irren",
    "fordern",
    "brauchen",
    "flehen",
    "boxen",
    "sprengen",
    "schlafen",
    "heben",
    "herrschen",
    "sterben",
    "geben",
    "zeigen",
    "dringen",
    "doppeln",
    "schmecken",
    "flüstern",
    "mauern",
    "schlagen",
    "kosten",
    "treffen",
    "füllen",
    "mögen",
    "lügen",
    "fischen",
    "schwächen",
    "wiegen",
    "fließen",
    "narren",
    "schwören",
    "klagen",
    "knien",
    "leeren",
    "schicken",
    "treiben",
    "haben",
    "warnen",
    "zeichnen",
    "führen",
    "schulden",
    "wandern",
    "prügeln",
    "wundern",
    "scheinen",
    "äußern",
    "sperren",
    "schenken",
    "lächeln",
    "zeugen",
    "grüßen",
    "wagen",
    "siegen",
    "altern",
    "kichern",
    "sinken",
    "fällen",
    "ärgern",
    "kürzen",
    "biegen",
    "helfen",
    "kehren",
    "rasen",
    "greifen",
    "gelten",
    "tauschen",
    "dauern",
    "jubeln",
    "trotzen",
    "beißen",
    "knarren",
    "dienen",
    "schreiten",
    "stammen",
    "zünden",
    "fragen",
    "spinnen",
    "läuten",
    "schrecken",
    "lockern",
    "bleiben",
    "graben",
    "fallen",
    "trennen",
    "rücken",
    "feiern",
    "zögern",
    "rufen",
    "runden",
    "planen",
    "ziehen",
    "münzen",
    "malen",
    "scheitern",
    "stecken",
    "bellen",
    "räumen",
    "reisen",
    "schützen",
    "erben",
    "werfen",
    "wüten",
    "fangen",
    "deuten",
    "saufen",
    "tollen",
    "sichern",
    "hindern",
    "backen",
    "zählen",
    "fahren",
    "kriegen",
    "werden"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}
